$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$values = @(
    @(1, 1, "27×79=2133"),
    @(1, 2, "56×38=2128"),
    @(1, 3, "28×84=2352"),
    @(1, 4, "65×59=3835"),
    @(1, 5, "87×32=2784"),
    @(5, 1, "11×98=1078"),
    @(5, 2, "99×26=2574"),
    @(5, 3, "24×49=1176"),
    @(5, 4, "28×68=1904"),
    @(5, 5, "14×91=1274"),
    @(10, 1, "29×64=1856"),
    @(10, 2, "67×53=3551"),
    @(10, 3, "30×16=480"),
    @(10, 4, "55×11=605"),
    @(10, 5, "68×83=5644"),
    @(15, 1, "24×11=264"),
    @(15, 2, "98×64=6272"),
    @(15, 3, "11×89=979"),
    @(15, 4, "52×16=832"),
    @(15, 5, "56×23=1288"),
    @(20, 1, "72×74=5328"),
    @(20, 2, "70×34=2380"),
    @(20, 3, "19×57=1083"),
    @(20, 4, "28×37=1036"),
    @(20, 5, "77×47=3619")
)

foreach ($entry in $values) {
    $row = $entry[0]
    $col = $entry[1]
    $newText = $entry[2]
    $cell = $t.Cell($row, $col)
    $cell.Range.Text = $newText
}
